# Update outdated "Uyumlu Marka" (Compatible Brand) values on the
# "Urun_Ozellik_Bilgileri" worksheet by appending " Uyumlu" to the brand
# names in column D (rows 2-13), matching the current Trendyol template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Urun_Ozellik_Bilgileri")

$ws.Range("D2").Value = "Amazfit Uyumlu"
$ws.Range("D3").Value = "Apple Uyumlu"
$ws.Range("D5").Value = "Fitbit Uyumlu"
$ws.Range("D6").Value = "Garmin Uyumlu"
$ws.Range("D7").Value = "Haylou Uyumlu"
$ws.Range("D8").Value = "Honor Uyumlu"
$ws.Range("D9").Value = "Huawei Uyumlu"
$ws.Range("D10").Value = "Oppo Uyumlu"
$ws.Range("D11").Value = "Samsung Uyumlu"
$ws.Range("D12").Value = "Winex Uyumlu"
$ws.Range("D13").Value = "Xiaomi Uyumlu"
